# Refactor synthetic array /3 for publipostage
# Replace the "black"/"noir" markers with "blue"/"bleu" markers across the sheet:
#   ⬛  -> 📘
#   🟧  -> 📙
#   noir -> bleu
# ("orange" stays untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $val = $cell.Value()

        # Only touch actual text cells (avoid PowerShell's bool<->string
        # coercion turning TRUE/FALSE cells into a match for any
        # non-empty string literal).
        if ($val -isnot [string]) {
            continue
        }

        if ($val -ceq "⬛") {
            $cell.Value = "📘"
        }
        elseif ($val -ceq "🟧") {
            $cell.Value = "📙"
        }
        elseif ($val -ceq "noir") {
            $cell.Value = "bleu"
        }
    }
}
